$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

Replace-Text "2024-10-22 Tuesday" "2024-10-23 Wednesday"

Replace-Text "398÷3=132, 2" "118÷5=23, 3"
Replace-Text "197÷5=39, 2" "859÷6=143, 1"
Replace-Text "715÷5=143, 0" "195÷6=32, 3"
Replace-Text "331÷5=66, 1" "307÷4=76, 3"
Replace-Text "980÷5=196, 0" "942÷4=235, 2"
Replace-Text "550÷2=275, 0" "782÷2=391, 0"
Replace-Text "801÷8=100, 1" "218÷4=54, 2"
Replace-Text "587÷3=195, 2" "541÷9=60, 1"
Replace-Text "976÷8=122, 0" "153÷8=19, 1"
Replace-Text "214÷3=71, 1" "684÷7=97, 5"
Replace-Text "540÷8=67, 4" "215÷9=23, 8"
Replace-Text "131÷3=43, 2" "314÷5=62, 4"
Replace-Text "824÷7=117, 5" "379÷9=42, 1"
Replace-Text "647÷8=80, 7" "834÷2=417, 0"
Replace-Text "229÷3=76, 1" "445÷6=74, 1"
Replace-Text "415÷9=46, 1" "449÷5=89, 4"
Replace-Text "285÷2=142, 1" "251÷9=27, 8"
Replace-Text "131÷5=26, 1" "834÷6=139, 0"
Replace-Text "714÷8=89, 2" "425÷3=141, 2"
Replace-Text "128÷6=21, 2" "102÷2=51, 0"
Replace-Text "877÷8=109, 5" "884÷7=126, 2"
Replace-Text "769÷3=256, 1" "126÷7=18, 0"
Replace-Text "992÷8=124, 0" "689÷2=344, 1"
Replace-Text "202÷5=40, 2" "887÷7=126, 5"
Replace-Text "383÷6=63, 5" "842÷6=140, 2"
